# Updates the SS-method comparison table: refreshed lead-time/accuracy
# metrics (now computed over 104 samples instead of 1641), reworked
# std_dev_ss (col K) so it is only populated for the "Rule-based SS"
# rows, and re-assigned Selected_method (col Q) per SKU/echelon. Also
# appends SKU_6 (DC / Store / WH) as rows 17-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: sku_id, echelon_type, lead_time, service_level, RMSE, MAE,
#          MAPE (%), Bias, WAPE (%), n_samples, std_dev_ss, Selected_method
# std_dev_ss is $null where the cell must end up blank (cleared).
$rows = @(
    @(2,  "SKU_1", "DC",    11, 0.9,  9.59,  7.88,               4.24,  1.64,   4.24,  104, $null,  "Forecast-based SS"),
    @(3,  "SKU_1", "Store", 6,  0.98, 5.96,  4.39,               4,     0.04,   4.08,  104, $null,  "Forecast-based SS"),
    @(4,  "SKU_1", "WH",    11, 0.95, 10.02, 8.140000000000001,  4.48,  -0.38,  4.58,  104, $null,  "Hybrid"),
    @(5,  "SKU_2", "DC",    8,  0.95, 37.66, 28.16,              12.15, -2.1,   12.18, 104, 267.61, "Rule-based SS"),
    @(6,  "SKU_2", "Store", 8,  0.9,  17.23, 14.22,              12.19, -1.07,  11.98, 104, $null,  "Hybrid"),
    @(7,  "SKU_2", "WH",    16, 0.98, 23.5,  18.18,              12.04, -1.64,  11.71, 104, 399.99, "Rule-based SS"),
    @(8,  "SKU_3", "DC",    13, 0.98, 40.64, 34.66,              20.18, 33.49,  20.31, 104, 299.73, "Rule-based SS"),
    @(9,  "SKU_3", "Store", 9,  0.95, 25.33, 21.86,              19.18, 21.27,  19.46, 104, 155.65, "Rule-based SS"),
    @(10, "SKU_3", "WH",    12, 0.9,  40.55, 34.14,              18.73, 32.5,   18.51, 104, 182.2,  "Rule-based SS"),
    @(11, "SKU_4", "DC",    17, 0.9,  52.81, 44.03,              20.89, -43.09, 20.58, 104, 247.78, "Rule-based SS"),
    @(12, "SKU_4", "Store", 15, 0.98, 23.14, 19.39,              19.94, -18.96, 19.66, 104, 175.16, "Rule-based SS"),
    @(13, "SKU_4", "WH",    9,  0.95, 38.29, 32.98,              18.56, -32.2,  18.2,  104, 232.5,  "Rule-based SS"),
    @(14, "SKU_5", "DC",    7,  0.95, 77.88, 60.31,              28.75, 1.96,   28.48, 104, 401.81, "Rule-based SS"),
    @(15, "SKU_5", "Store", 5,  0.9,  46.1,  36.15,              27.64, 1.12,   27.81, 104, 169.75, "Rule-based SS"),
    @(16, "SKU_5", "WH",    8,  0.98, 65.54000000000001, 51.7,   26.73, -11.27, 27.02, 104, 404.05, "Rule-based SS"),
    @(17, "SKU_6", "DC",    5,  0.98, 14.19, 11.45,              4.65,  -0.58,  4.66,  104, $null,  "Hybrid"),
    @(18, "SKU_6", "Store", 5,  0.95, 4.85,  3.82,               3.66,  0.03,   3.56,  104, $null,  "Forecast-based SS"),
    @(19, "SKU_6", "WH",    8,  0.9,  7.41,  5.66,               3.42,  -0.38,  3.46,  104, $null,  "Forecast-based SS")
)

foreach ($row in $rows) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value  = $row[1]   # sku_id
    $ws.Cells.Item($r, 2).Value  = $row[2]   # echelon_type
    $ws.Cells.Item($r, 3).Value  = $row[3]   # lead_time
    $ws.Cells.Item($r, 4).Value  = $row[4]   # service_level
    $ws.Cells.Item($r, 5).Value  = $row[5]   # RMSE
    $ws.Cells.Item($r, 6).Value  = $row[6]   # MAE
    $ws.Cells.Item($r, 7).Value  = $row[7]   # MAPE (%)
    $ws.Cells.Item($r, 8).Value  = $row[8]   # Bias
    $ws.Cells.Item($r, 9).Value  = $row[9]   # WAPE (%)
    $ws.Cells.Item($r, 10).Value = $row[10]  # n_samples

    $stdDev = $row[11]
    if ($null -eq $stdDev) {
        $ws.Cells.Item($r, 11).ClearContents()  # std_dev_ss
    } else {
        $ws.Cells.Item($r, 11).Value = $stdDev
    }

    # rmse_ss, mae_ss, hybrid_with_no_var_ss, hybrid_with_var_ss, ml_ss stay blank
    $ws.Cells.Item($r, 12).ClearContents()
    $ws.Cells.Item($r, 13).ClearContents()
    $ws.Cells.Item($r, 14).ClearContents()
    $ws.Cells.Item($r, 15).ClearContents()
    $ws.Cells.Item($r, 16).ClearContents()

    $ws.Cells.Item($r, 17).Value = $row[12]  # Selected_method
}
